# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps for the
# a4807e53-... row (row 4) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-18 12:17:24"
$wsZhCn.Range("G4").Value = "2016-01-18 12:18:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-18 12:17:33"
$wsDeDe.Range("G4").Value = "2016-01-18 12:18:25"
